$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = [double]"3"
$ws.Range("F2").Value = [double]"1"
$ws.Range("G2").Value = [double]"1.172733"
$ws.Range("H2").Value = [double]"3.518199"
$ws.Range("I2").Value = [double]"0.02094912533655296"
$ws.Range("J2").Value = [double]"0.02094912533655295"
$ws.Range("M2").Value = [double]"2.157506"
$ws.Range("N2").Value = [double]"6.472517999999999"
$ws.Range("O2").Value = [double]"0.3549648016839517"
$ws.Range("P2").Value = [double]"0.3549648016839516"
$ws.Range("Q2").Value = [double]"2.530178483898"
$ws.Range("R2").Value = [double]"22.771606355082"
$ws.Range("S2").Value = [double]"0.007436202120541768"
$ws.Range("T2").Value = [double]"0.007436202120541766"
$ws.Range("E3").Value = [double]"3"
$ws.Range("F3").Value = [double]"1"
$ws.Range("G3").Value = [double]"1.172733"
$ws.Range("H3").Value = [double]"3.518199"
$ws.Range("I3").Value = [double]"0.02094912533655296"
$ws.Range("J3").Value = [double]"0.02094912533655295"
$ws.Range("O3").Value = [double]"0.4793705560628122"
$ws.Range("P3").Value = [double]"0.4793705560628121"
$ws.Range("Q3").Value = [double]"3.416938978204"
$ws.Range("R3").Value = [double]"30.752450803836"
$ws.Range("S3").Value = [double]"0.01004239386161294"
$ws.Range("T3").Value = [double]"0.01004239386161294"
$ws.Range("E4").Value = [double]"3"
$ws.Range("F4").Value = [double]"1"
$ws.Range("G4").Value = [double]"1.172733"
$ws.Range("H4").Value = [double]"3.518199"
$ws.Range("I4").Value = [double]"0.02094912533655296"
$ws.Range("J4").Value = [double]"0.02094912533655295"
$ws.Range("M4").Value = [double]"0.018986"
$ws.Range("N4").Value = [double]"0.05695799999999999"
$ws.Range("O4").Value = [double]"0.003123681567871193"
$ws.Range("P4").Value = [double]"0.003123681567871192"
$ws.Range("Q4").Value = [double]"0.022265508738"
$ws.Range("R4").Value = [double]"0.200389578642"
$ws.Range("S4").Value = [double]"6.543839667681388E-05"
$ws.Range("T4").Value = [double]"6.543839667681385E-05"
$ws.Range("E5").Value = [double]"3"
$ws.Range("F5").Value = [double]"1"
$ws.Range("G5").Value = [double]"1.172733"
$ws.Range("H5").Value = [double]"3.518199"
$ws.Range("I5").Value = [double]"0.02094912533655296"
$ws.Range("J5").Value = [double]"0.02094912533655295"
$ws.Range("M5").Value = [double]"0.9848966666666668"
$ws.Range("N5").Value = [double]"2.95469"
$ws.Range("O5").Value = [double]"0.1620406385718132"
$ws.Range("P5").Value = [double]"0.1620406385718132"
$ws.Range("Q5").Value = [double]"1.15502082259"
$ws.Range("R5").Value = [double]"10.39518740331"
$ws.Range("S5").Value = [double]"0.003394609647055993"
$ws.Range("T5").Value = [double]"0.003394609647055992"
$ws.Range("D6").Value = "Resolving-Mac"
$ws.Range("E6").Value = [double]"3"
$ws.Range("F6").Value = [double]"1"
$ws.Range("G6").Value = [double]"1.172733"
$ws.Range("H6").Value = [double]"3.518199"
$ws.Range("I6").Value = [double]"0.02094912533655296"
$ws.Range("J6").Value = [double]"0.02094912533655295"
$ws.Range("M6").Value = [double]"0.003041"
$ws.Range("N6").Value = [double]"0.009122999999999999"
$ws.Range("O6").Value = [double]"0.0005003221135518961"
$ws.Range("P6").Value = [double]"0.000500322113551896"
$ws.Range("Q6").Value = [double]"0.003566281053"
$ws.Range("R6").Value = [double]"0.032096529477"
$ws.Range("S6").Value = [double]"1.048131066544775E-05"
$ws.Range("T6").Value = [double]"1.048131066544775E-05"
$ws.Range("G7").Value = [double]"54.58029933333334"
$ws.Range("H7").Value = [double]"163.740898"
$ws.Range("I7").Value = [double]"0.974995614211059"
$ws.Range("J7").Value = [double]"0.974995614211059"
$ws.Range("M7").Value = [double]"2.157506"
$ws.Range("N7").Value = [double]"6.472517999999999"
$ws.Range("O7").Value = [double]"0.3549648016839517"
$ws.Range("P7").Value = [double]"0.3549648016839516"
$ws.Range("Q7").Value = [double]"117.7573232934627"
$ws.Range("R7").Value = [double]"1059.815909641164"
$ws.Range("S7").Value = [double]"0.3460891248411512"
$ws.Range("T7").Value = [double]"0.3460891248411511"
$ws.Range("G8").Value = [double]"54.58029933333334"
$ws.Range("H8").Value = [double]"163.740898"
$ws.Range("I8").Value = [double]"0.974995614211059"
$ws.Range("J8").Value = [double]"0.974995614211059"
$ws.Range("O8").Value = [double]"0.4793705560628122"
$ws.Range("P8").Value = [double]"0.4793705560628121"
$ws.Range("Q8").Value = [double]"159.0281438606302"
$ws.Range("R8").Value = [double]"1431.253294745672"
$ws.Range("S8").Value = [double]"0.4673841897431584"
$ws.Range("T8").Value = [double]"0.4673841897431584"
$ws.Range("G9").Value = [double]"54.58029933333334"
$ws.Range("H9").Value = [double]"163.740898"
$ws.Range("I9").Value = [double]"0.974995614211059"
$ws.Range("J9").Value = [double]"0.974995614211059"
$ws.Range("M9").Value = [double]"0.018986"
$ws.Range("N9").Value = [double]"0.05695799999999999"
$ws.Range("O9").Value = [double]"0.003123681567871193"
$ws.Range("P9").Value = [double]"0.003123681567871192"
$ws.Range("Q9").Value = [double]"1.036261563142667"
$ws.Range("R9").Value = [double]"9.326354068283999"
$ws.Range("S9").Value = [double]"0.003045575828866338"
$ws.Range("T9").Value = [double]"0.003045575828866337"
$ws.Range("G10").Value = [double]"54.58029933333334"
$ws.Range("H10").Value = [double]"163.740898"
$ws.Range("I10").Value = [double]"0.974995614211059"
$ws.Range("J10").Value = [double]"0.974995614211059"
$ws.Range("M10").Value = [double]"0.9848966666666668"
$ws.Range("N10").Value = [double]"2.95469"
$ws.Range("O10").Value = [double]"0.1620406385718132"
$ws.Range("P10").Value = [double]"0.1620406385718132"
$ws.Range("Q10").Value = [double]"53.7559548790689"
$ws.Range("R10").Value = [double]"483.8035939116201"
$ws.Range("S10").Value = [double]"0.1579889119314772"
$ws.Range("T10").Value = [double]"0.1579889119314772"
$ws.Range("D11").Value = "Resolving-Mac"
$ws.Range("G11").Value = [double]"54.58029933333334"
$ws.Range("H11").Value = [double]"163.740898"
$ws.Range("I11").Value = [double]"0.974995614211059"
$ws.Range("J11").Value = [double]"0.974995614211059"
$ws.Range("M11").Value = [double]"0.003041"
$ws.Range("N11").Value = [double]"0.009122999999999999"
$ws.Range("O11").Value = [double]"0.0005003221135518961"
$ws.Range("P11").Value = [double]"0.000500322113551896"
$ws.Range("Q11").Value = [double]"0.1659786902726667"
$ws.Range("R11").Value = [double]"1.493808212454"
$ws.Range("S11").Value = [double]"0.0004878118664059061"
$ws.Range("T11").Value = [double]"0.000487811866405906"
$ws.Range("E12").Value = [double]"1"
$ws.Range("F12").Value = [double]"0.3333333333333333"
$ws.Range("G12").Value = [double]"0.2270136666666667"
$ws.Range("H12").Value = [double]"0.681041"
$ws.Range("I12").Value = [double]"0.0040552604523881"
$ws.Range("J12").Value = [double]"0.0040552604523881"
$ws.Range("M12").Value = [double]"2.157506"
$ws.Range("N12").Value = [double]"6.472517999999999"
$ws.Range("O12").Value = [double]"0.3549648016839517"
$ws.Range("P12").Value = [double]"0.3549648016839516"
$ws.Range("Q12").Value = [double]"0.4897833479153333"
$ws.Range("R12").Value = [double]"4.408050131237999"
$ws.Range("S12").Value = [double]"0.001439474722258714"
$ws.Range("T12").Value = [double]"0.001439474722258714"
$ws.Range("E13").Value = [double]"1"
$ws.Range("F13").Value = [double]"0.3333333333333333"
$ws.Range("G13").Value = [double]"0.2270136666666667"
$ws.Range("H13").Value = [double]"0.681041"
$ws.Range("I13").Value = [double]"0.0040552604523881"
$ws.Range("J13").Value = [double]"0.0040552604523881"
$ws.Range("O13").Value = [double]"0.4793705560628122"
$ws.Range("P13").Value = [double]"0.4793705560628121"
$ws.Range("Q13").Value = [double]"0.6614394292804444"
$ws.Range("R13").Value = [double]"5.952954863524"
$ws.Range("S13").Value = [double]"0.001943972458040815"
$ws.Range("T13").Value = [double]"0.001943972458040815"
$ws.Range("E14").Value = [double]"1"
$ws.Range("F14").Value = [double]"0.3333333333333333"
$ws.Range("G14").Value = [double]"0.2270136666666667"
$ws.Range("H14").Value = [double]"0.681041"
$ws.Range("I14").Value = [double]"0.0040552604523881"
$ws.Range("J14").Value = [double]"0.0040552604523881"
$ws.Range("M14").Value = [double]"0.018986"
$ws.Range("N14").Value = [double]"0.05695799999999999"
$ws.Range("O14").Value = [double]"0.003123681567871193"
$ws.Range("P14").Value = [double]"0.003123681567871192"
$ws.Range("Q14").Value = [double]"0.004310081475333333"
$ws.Range("R14").Value = [double]"0.038790733278"
$ws.Range("S14").Value = [double]"1.26673423280417E-05"
$ws.Range("T14").Value = [double]"1.26673423280417E-05"
$ws.Range("E15").Value = [double]"1"
$ws.Range("F15").Value = [double]"0.3333333333333333"
$ws.Range("G15").Value = [double]"0.2270136666666667"
$ws.Range("H15").Value = [double]"0.681041"
$ws.Range("I15").Value = [double]"0.0040552604523881"
$ws.Range("J15").Value = [double]"0.0040552604523881"
$ws.Range("M15").Value = [double]"0.9848966666666668"
$ws.Range("N15").Value = [double]"2.95469"
$ws.Range("O15").Value = [double]"0.1620406385718132"
$ws.Range("P15").Value = [double]"0.1620406385718132"
$ws.Range("Q15").Value = [double]"0.2235850035877778"
$ws.Range("R15").Value = [double]"2.01226503229"
$ws.Range("S15").Value = [double]"0.0006571169932799879"
$ws.Range("T15").Value = [double]"0.0006571169932799878"
$ws.Range("D16").Value = "Resolving-Mac"
$ws.Range("E16").Value = [double]"1"
$ws.Range("F16").Value = [double]"0.3333333333333333"
$ws.Range("G16").Value = [double]"0.2270136666666667"
$ws.Range("H16").Value = [double]"0.681041"
$ws.Range("I16").Value = [double]"0.0040552604523881"
$ws.Range("J16").Value = [double]"0.0040552604523881"
$ws.Range("M16").Value = [double]"0.003041"
$ws.Range("N16").Value = [double]"0.009122999999999999"
$ws.Range("O16").Value = [double]"0.0005003221135518961"
$ws.Range("P16").Value = [double]"0.000500322113551896"
$ws.Range("Q16").Value = [double]"0.0006903485603333334"
$ws.Range("R16").Value = [double]"0.006213137043"
$ws.Range("S16").Value = [double]"2.028936480542232E-06"
$ws.Range("T16").Value = [double]"2.028936480542232E-06"
